$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(162, 1).Value = "2023-12-10 14:17:24"
$ws.Cells.Item(162, 2).Value = 0.0004

$ws.Cells.Item(163, 1).Value = "2023-12-10 14:17:36"
$ws.Cells.Item(163, 2).Value = 0.0004
